$wb = $excel.ActiveWorkbook

# Update the value on the "DatosCuenta" sheet (D2: 120 -> 121)
$wsCuenta = $wb.Worksheets.Item("DatosCuenta")
$wsCuenta.Range("D2").Value = 121

# Make "DatosCuenta" the active sheet/tab, with D3 selected
$wsCuenta.Activate()
$wsCuenta.Range("D3").Select()
